$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("B2").Value = 4
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 1

# Remove the now-obsolete row 5 (shifts nothing below it, just clears it)
$ws.Range("A5:B5").EntireRow.Delete()
